# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.877.58'
$ws.Range("E2").Value = '  -5.60%  '
$ws.Range("D3").Value = '2.980.30'
$ws.Range("E3").Value = '  -6.09%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '124.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.43%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").Value = '2.972.02'
$ws.Range("E8").Value = '  -6.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("E10").Value = '  -8.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.75%  '
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000219'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.40%  '
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '3.471.50'
$ws.Range("E16").Value = '  -6.07%  '
$ws.Range("D17").Value = '2.973.02'
$ws.Range("E17").Value = '  -6.21%  '
$ws.Range("D18").Value = '59.828.22'
$ws.Range("E18").Value = '  -5.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '425.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.667'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.89%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.14%  '
$ws.Range("E29").Value = '  -7.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0965'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.917'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -18.61%  '
$ws.Range("E38").Value = '  +3.83%  '
$ws.Range("D39").Value = '0.0₃0647'
$ws.Range("E39").Value = '  -11.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0353'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.46%  '
$ws.Range("E41").Value = '  -5.94%  '
$ws.Range("D42").Value = '2.658.31'
$ws.Range("E42").Value = '  -5.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '366.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.93%  '
$ws.Range("E44").Value = '  -8.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '120.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.232'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.35%  '
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.09%  '
$ws.Range("E51").Value = '  -8.12%  '
